# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates to match the target diff across
# sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW (WVR unchanged).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1090.4884
$ws.Range("I15").Value = 1090.4884
$ws.Range("K15").Value = 3271.4652
$ws.Range("M15").Value = -3102.4652

$ws.Range("H38").Value = 1024.5714
$ws.Range("J38").Value = 3000
$ws.Range("L38").Value = 9000
$ws.Range("N38").Value = -9744

$ws.Range("H40").Value = 2437.5
$ws.Range("J40").Value = 1875
$ws.Range("L40").Value = 1875
$ws.Range("N40").Value = -2225

$ws.Range("H51").Value = 6600.2
$ws.Range("I51").Value = 5001
$ws.Range("K51").Value = 5001
$ws.Range("M51").Value = -4517

$ws.Range("H58").Value = 1105.4667
$ws.Range("I58").Value = 455.42856
$ws.Range("J58").Value = 1674.25
$ws.Range("K58").Value = 1366.28568
$ws.Range("L58").Value = 5022.75
$ws.Range("M58").Value = -1216.28568
$ws.Range("N58").Value = -5322.75

$ws.Range("H70").Value = 13444
$ws.Range("J70").Value = 15264.571
$ws.Range("L70").Value = 45793.713
$ws.Range("N70").Value = -46333.713

$ws.Range("H73").Value = 13444
$ws.Range("J73").Value = 15264.571
$ws.Range("L73").Value = 45793.713
$ws.Range("N73").Value = -47665.713

$ws.Range("H86").Value = 1837
$ws.Range("J86").Value = 1874
$ws.Range("L86").Value = 1874
$ws.Range("N86").Value = -4120

$ws.Range("H89").Value = 1837
$ws.Range("J89").Value = 1874
$ws.Range("L89").Value = 9370
$ws.Range("N89").Value = -20602

$ws.Range("H121").Value = 150
$ws.Range("I121").Value = 150
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 450
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 1297
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2842.896
$ws.Range("I32").Value = 2050.147
$ws.Range("K32").Value = 2050.147
$ws.Range("M32").Value = -1763.147

$ws.Range("H74").Value = 1739.3096
$ws.Range("I74").Value = 1637.8948
$ws.Range("K74").Value = 1637.8948
$ws.Range("M74").Value = -763.8948

$ws.Range("H77").Value = 1739.3096
$ws.Range("I77").Value = 1637.8948
$ws.Range("K77").Value = 8189.474
$ws.Range("M77").Value = -3821.474

$ws.Range("H132").Value = 2169.5925
$ws.Range("I132").Value = 1749.5
$ws.Range("K132").Value = 5248.5
$ws.Range("M132").Value = -2718.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 33700
$ws.Range("I82").Value = 12500
$ws.Range("K82").Value = 12500
$ws.Range("M82").Value = -12117

$ws.Range("H85").Value = 33700
$ws.Range("I85").Value = 12500
$ws.Range("K85").Value = 12500
$ws.Range("M85").Value = -11174

$ws.Range("H86").Value = 102353.7
$ws.Range("I86").Value = 2344.3845
$ws.Range("J86").Value = 288085.28
$ws.Range("K86").Value = 2344.3845
$ws.Range("L86").Value = 288085.28
$ws.Range("M86").Value = -1221.3845
$ws.Range("N86").Value = -290331.28

$ws.Range("H89").Value = 102353.7
$ws.Range("I89").Value = 2344.3845
$ws.Range("J89").Value = 288085.28
$ws.Range("K89").Value = 11721.9225
$ws.Range("L89").Value = 1440426.4
$ws.Range("M89").Value = -6105.922500000001
$ws.Range("N89").Value = -1451658.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1933.36
$ws.Range("I31").Value = 1067.9286
$ws.Range("J31").Value = 3034.818
$ws.Range("K31").Value = 1067.9286
$ws.Range("L31").Value = 3034.818
$ws.Range("M31").Value = -772.9286
$ws.Range("N31").Value = -3624.818

$ws.Range("H34").Value = 1933.36
$ws.Range("I34").Value = 1067.9286
$ws.Range("J34").Value = 3034.818
$ws.Range("K34").Value = 1067.9286
$ws.Range("L34").Value = 3034.818
$ws.Range("M34").Value = -865.9286
$ws.Range("N34").Value = -3438.818

$ws.Range("H86").Value = 1998.4546
$ws.Range("I86").Value = 1577.6
$ws.Range("J86").Value = 2349.1667
$ws.Range("K86").Value = 1577.6
$ws.Range("L86").Value = 2349.1667
$ws.Range("M86").Value = -454.5999999999999
$ws.Range("N86").Value = -4595.1667

$ws.Range("H89").Value = 1998.4546
$ws.Range("I89").Value = 1577.6
$ws.Range("J89").Value = 2349.1667
$ws.Range("K89").Value = 7888
$ws.Range("L89").Value = 11745.8335
$ws.Range("M89").Value = -2272
$ws.Range("N89").Value = -22977.8335

$ws.Range("H122").Value = 1702
$ws.Range("I122").Value = 1604
$ws.Range("J122").Value = 2633
$ws.Range("K122").Value = 4812
$ws.Range("L122").Value = 7899
$ws.Range("M122").Value = -2362
$ws.Range("N122").Value = -12799

$ws.Range("H134").Value = 870.75
$ws.Range("I134").Value = 870.75
$ws.Range("K134").Value = 2612.25
$ws.Range("M134").Value = -77.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 362.75
$ws.Range("J23").Value = 362.75
$ws.Range("L23").Value = 1088.25
$ws.Range("N23").Value = -1558.25

$ws.Range("H48").Value = 1499.6
$ws.Range("J48").Value = 1499.6
$ws.Range("L48").Value = 4498.799999999999
$ws.Range("N48").Value = -4998.799999999999

$ws.Range("H131").Value = 770.09
$ws.Range("I131").Value = 515.6667
$ws.Range("J131").Value = 786.3298
$ws.Range("K131").Value = 1547.0001
$ws.Range("L131").Value = 2358.9894
$ws.Range("M131").Value = 3492.9999
$ws.Range("N131").Value = -12438.9894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7695906.5
$ws.Range("I132").Value = 19232768
$ws.Range("J132").Value = 4665.6665
$ws.Range("K132").Value = 57698304
$ws.Range("L132").Value = 13996.9995
$ws.Range("M132").Value = -57695774
$ws.Range("N132").Value = -19056.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2621.7
$ws.Range("I82").Value = 1748.8
$ws.Range("K82").Value = 1748.8
$ws.Range("M82").Value = -1387.8

$ws.Range("H85").Value = 2621.7
$ws.Range("I85").Value = 1748.8
$ws.Range("K85").Value = 1748.8
$ws.Range("M85").Value = -500.8

$ws.Range("H132").Value = 1889.7576
$ws.Range("I132").Value = 1848.5454
$ws.Range("J132").Value = 1910.3636
$ws.Range("K132").Value = 5545.6362
$ws.Range("L132").Value = 5731.0908
$ws.Range("M132").Value = -3015.6362
$ws.Range("N132").Value = -10791.0908

$ws.Range("H136").Value = 2768.9092
$ws.Range("J136").Value = 4432.8887
$ws.Range("L136").Value = 13298.6661
$ws.Range("N136").Value = -18398.6661
